$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B slightly to fit new content
# (ColumnWidth quantizes internally to sixths of a character; 48.6666...
# is the closest settable value to the target stored width of 49.42578125)
$ws.Columns.Item(2).ColumnWidth = 48.666666666666664

# Add the two new rows of data
$ws.Range("B16").Value = "Pridejau telefono numeri su ngx mask"
$ws.Range("C16").Value = 1

$ws.Range("B17").Value = "Padariau kad isikeltu daug paveiksleliu ir isirasytu i db"
$ws.Range("C17").Value = 2
